$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: values are written as plain (non-exponential) decimal strings so that
# Excel stores the exact underlying double without re-tagging the cells with a
# "Scientific" number format (which would otherwise happen if an "E" notation
# string were assigned directly).

$ws.Range("D2").Value = "0.00000000004077478427296493"
$ws.Range("E2").Value = "0.00000000004077478427296493"

$ws.Range("D3").Value = "0.0000000000000000000000000000000000000000000000000000000001800806493465486"
$ws.Range("E3").Value = "0.0000000000000000000000000000000000000000000000000000000001800806493465486"

$ws.Range("D4").Value = "0.0000000000000000000000000003612050650184042"
$ws.Range("E4").Value = "0.0000000000000000000000000003612050650184042"

$ws.Range("D5").Value = "0.0000000000000000076179771115205"
$ws.Range("E5").Value = "0.0000000000000000076179771115205"

$ws.Range("D6").Value = "0.9999999997137472"
$ws.Range("E6").Value = "0.9999999997137472"

$ws.Range("D8").Value = "0.9999999997451001"
$ws.Range("E8").Value = "0.0000000002548998789819734"

$ws.Range("D10").Value = "0.00000000000000000326597792734769"
$ws.Range("E10").Value = "1"

$ws.Range("D11").Value = "0.000000000000000000000000000000000000000000000000000000000000000000000000000000000000000001268111908755485"

$ws.Range("F11").Value = "26.9232234954834"
